# Insert two new rows at row 170, shifting existing rows 170:243 down to 172:245
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(170).Resize(2).Insert()

# New row 170: Camote, 1a (cosecha), Region del Maule
$ws.Cells.Item(170, 1).Value = 7
$ws.Cells.Item(170, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(170, 3).Value = "Ñuble"
$ws.Cells.Item(170, 4).Value = 44992
$ws.Cells.Item(170, 5).Value = 16
$ws.Cells.Item(170, 6).Value = 100112045
$ws.Cells.Item(170, 7).Value = "Zapallo"
$ws.Cells.Item(170, 8).Value = "Camote"
$ws.Cells.Item(170, 9).Value = "1a (cosecha)"
$ws.Cells.Item(170, 10).Value = 300
$ws.Cells.Item(170, 11).Value = 450
$ws.Cells.Item(170, 12).Value = 500
$ws.Cells.Item(170, 13).Value = 475
$ws.Cells.Item(170, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(170, 15).Value = "Región del Maule"
$ws.Cells.Item(170, 16).Value = 475
$ws.Cells.Item(170, 17).Value = 1
$ws.Cells.Item(170, 18).Value = "Hortaliza"

# New row 171: Paine, 1a (cosecha), Region del Maule
$ws.Cells.Item(171, 1).Value = 7
$ws.Cells.Item(171, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(171, 3).Value = "Ñuble"
$ws.Cells.Item(171, 4).Value = 44992
$ws.Cells.Item(171, 5).Value = 16
$ws.Cells.Item(171, 6).Value = 100112045
$ws.Cells.Item(171, 7).Value = "Zapallo"
$ws.Cells.Item(171, 8).Value = "Paine"
$ws.Cells.Item(171, 9).Value = "1a (cosecha)"
$ws.Cells.Item(171, 10).Value = 150
$ws.Cells.Item(171, 11).Value = 350
$ws.Cells.Item(171, 12).Value = 350
$ws.Cells.Item(171, 13).Value = 350
$ws.Cells.Item(171, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(171, 15).Value = "Región del Maule"
$ws.Cells.Item(171, 16).Value = 350
$ws.Cells.Item(171, 17).Value = 1
$ws.Cells.Item(171, 18).Value = "Hortaliza"
